# Publication prep for release 0.2.0 (CodeSystem-eclaire-study-phase-code-system)
#
# Changes on the "Metadata" sheet:
#   - Version bumped: 0.1.1 -> 0.2.0
#   - Date bumped:    2023-10-19T16:17:18+00:00 -> 2023-10-19T17:05:12+00:00
#   - New "Jurisdiction" / "iso:code:3166:FR" row inserted right after the
#     "Contact" row (pushes Description..Count down by one row)
#
# The "Concepts" sheet is untouched content-wise.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Version ---------------------------------------------------------------
$ws.Cells.Item(3, 2).Value = "0.2.0"

# --- Date --------------------------------------------------------------
$ws.Cells.Item(8, 2).Value = "2023-10-19T17:05:12+00:00"

# --- Insert "Jurisdiction" row after "Contact" (row 10) --------------------
$ws.Rows.Item(11).Insert()

# Copy the formatting of the row above (Contact) onto the freshly inserted
# row so it keeps the same style (border/alignment) as the rest of the table.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = "iso:code:3166:FR"
